# Generate Report for Handoff
#
# The two previously-pending files (the .png dependency pair that was
# "IsDependency" on a725933f-....md) have completed handoff/handback and
# are now reported under their own new GUID names:
#   a725933f-ffae-4462-b8a5-da06c4540a1d.md  -> 3ca43955-bd38-4955-aafa-69ccbe78ff67.md
#   (replaces the old "b003c385-....png")    -> bbeb3347-5d2c-450a-b41a-f9167885289b.md
# The two *.png rows (which only existed because they were dependencies)
# drop out of the report entirely, and ".localization-config" shifts up
# to take the now-freed last row.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Overview"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop the now-empty trailing row (.localization-config moves up to row 4).
$ws1.Rows.Item(5).Delete()

# Row 2 / Row 3: same status text, just the backing file renamed.
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

# Row 4 now holds what used to be row 5's content.
$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

# Rebuild the hyperlinks: same relationship targets (positionally reused),
# new display text, one fewer link than before.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/3dd73262-558d-431e-a720-20d985d050e7.png", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/a725933f-ffae-4462-b8a5-da06c4540a1d.md", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/b003c385-feec-4394-98cc-eee6c26800de.png", "", "", ".localization-config") | Out-Null

# ----------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(5).Delete()

$ws2.Range("D2").Value = "2016-03-08 10:42:36"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"
$ws2.Range("I2").ClearContents()

$ws2.Range("D3").Value = "2016-03-08 10:42:36"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("C4").ClearContents()
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"
$ws2.Range("I4").ClearContents()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/3dd73262-558d-431e-a720-20d985d050e7.png", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20ac2adba5604147130bc10d36d398c41a5f34d1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d29f000ac868c69ef581cedf711e45f6d0e2b1c4.png", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.8b726d99015d6d04615f5d4a555d51fa6ca19b07.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/a725933f-ffae-4462-b8a5-da06c4540a1d.md", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20ac2adba5604147130bc10d36d398c41a5f34d1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a725933f-ffae-4462-b8a5-da06c4540a1d.21655dfe974d50af8ddbfb8a71e952d4a935ced3.zh-cn.xlf", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.1a9e9f6d8ea9306613957ff624dc3b666edb2bf1.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/b003c385-feec-4394-98cc-eee6c26800de.png", "", "", ".localization-config") | Out-Null

# ----------------------------------------------------------------------
# Sheet 3: "de-de"
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(5).Delete()

$ws3.Range("D2").Value = "2016-03-08 10:42:40"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"
$ws3.Range("I2").ClearContents()

$ws3.Range("D3").Value = "2016-03-08 10:42:40"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("C4").ClearContents()
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"
$ws3.Range("I4").ClearContents()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/3dd73262-558d-431e-a720-20d985d050e7.png", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86dd9ae796bb7a71b95d23cd5091db3055352f6f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d29f000ac868c69ef581cedf711e45f6d0e2b1c4.png", "", "", "3ca43955-bd38-4955-aafa-69ccbe78ff67.8b726d99015d6d04615f5d4a555d51fa6ca19b07.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/a725933f-ffae-4462-b8a5-da06c4540a1d.md", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86dd9ae796bb7a71b95d23cd5091db3055352f6f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a725933f-ffae-4462-b8a5-da06c4540a1d.21655dfe974d50af8ddbfb8a71e952d4a935ced3.de-de.xlf", "", "", "bbeb3347-5d2c-450a-b41a-f9167885289b.1a9e9f6d8ea9306613957ff624dc3b666edb2bf1.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f23d81cec00fdbcb1361845c19daec5cb639910e/e2e/b003c385-feec-4394-98cc-eee6c26800de.png", "", "", ".localization-config") | Out-Null

Write-Output "Report for handoff generated."
